# Update the stale AMZN open/close price quote (row 2, columns F/G) and
# leave the active selection on the cell the user last edited (G2 -
# close_price), matching the saved worksheet view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# open_price (F2) and close_price (G2) refreshed to current values
$ws.Range("F2").Value = 147
$ws.Range("G2").Value = 149.75

# Move/save the selection on G2, as recorded in the sheet view
$ws.Range("G2").Select()
